$d = $word.ActiveDocument
$script:ridCounter = 9

function Add-EmptyPara() {
    $lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastP.Range.InsertParagraphAfter()
}

function Add-TextPara([string]$text, [bool]$bold) {
    $lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastP.Range.InsertParagraphAfter()
    $n = $d.Paragraphs.Count
    $newp = $d.Paragraphs.Item($n)
    $r = $d.Range($newp.Range.Start, $newp.Range.Start)
    $r.InsertAfter($text)
    $pp = $d.Paragraphs.Item($d.Paragraphs.Count)
    if ($bold) {
        $pp.Range.Font.Bold = 1
    } else {
        $pp.Range.Font.Bold = 0
    }
}

function Add-HyperlinkPara([string]$url) {
    $lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastP.Range.InsertParagraphAfter()
    $n = $d.Paragraphs.Count
    $newp = $d.Paragraphs.Item($n)
    $r0 = $d.Range($newp.Range.Start, $newp.Range.Start)
    $h = $d.Hyperlinks.Add($r0, $url, "", "", $url)
    $hh = $d.Hyperlinks.Item($d.Hyperlinks.Count)
    $rspace = $d.Range($hh.Range.End, $hh.Range.End)
    $rspace.InsertAfter(" ")

    # Rebuild the paragraph's content with the exact target formatting
    # (direct color/underline rather than a Hyperlink style, lower-case hex,
    # and no leftover placeholder run) while keeping the same relationship id
    # that Hyperlinks.Add just minted for this paragraph.
    $para = $d.Paragraphs.Item($d.Paragraphs.Count)
    $rid = "rId" + $script:ridCounter
    $script:ridCounter = $script:ridCounter + 1
    $full = $d.Range($para.Range.Start, $para.Range.End - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:rPr/></w:pPr><w:hyperlink r:id="' + $rid + '"><w:r><w:rPr><w:color w:val="1155cc"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $url + '</w:t></w:r></w:hyperlink><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($xml)
}

function Add-FdBlock([string]$label, [string]$startTime, [string]$url) {
    Add-TextPara $label $true
    Add-TextPara "Topic: [AND - 12] Chapter 5 - Utsman" $false
    Add-TextPara $startTime $false
    Add-EmptyPara
    Add-TextPara "Meeting Recording:" $false
    Add-HyperlinkPara $url
}

# Divider paragraph between the existing FD 3 block and the new FD 4 block
Add-EmptyPara

Add-FdBlock "FD 4" "Start Time: Apr 12, 2022 07:03 PM" "https://us02web.zoom.us/rec/share/Rt3so6Kf3Z4mMLPYdzv4ukajelJMMPutf1_KDrj6MIHDV6XhURLV20FI5-E3JtAQ.o5q5VLd_JaQVu8gg"

Add-EmptyPara

Add-FdBlock "FD 5" "Start Time: Apr 13, 2022 07:03 PM" "https://us02web.zoom.us/rec/share/3CnrlClKXHEByi6Pzt44ABRdp7BJ-kLc5W9rhMIuCNhtk2Rkx_4Mzm22j9BhLKhm.th3JB-V4GCmM-HPU"

Add-EmptyPara

Add-FdBlock "FD 6" "Start Time: Apr 14, 2022 07:04 PM" "https://us02web.zoom.us/rec/share/lIwDCuldyMgcwDGICqd4hPJGj4qZ2VrY48aWOT5MywGloutdi5nN0w5A-JwoSHU_.bHz3mjV-XEjsgKuh"

Write-Output "Done. Paragraphs: $($d.Paragraphs.Count)"
